$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2077
$ws1.Range("F6").Value = 622
$ws1.Range("F9").Value = 10668
$ws1.Range("F12").Value = 285
$ws1.Range("F15").Value = 7530

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2077
$ws4.Range("F6").Value = 622
$ws4.Range("F12").Value = 10668
$ws4.Range("F15").Value = 285
$ws4.Range("F18").Value = 7530
